$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header row: "_old" -> "_FV2310", "_new" -> "_FV2404"
for ($c = 1; $c -le 21; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $v = $cell.Value2
    if ($v -ne $null) {
        if ($v.EndsWith("_old")) {
            $cell.Value = $v.Substring(0, $v.Length - 4) + "_FV2310"
        } elseif ($v.EndsWith("_new")) {
            $cell.Value = $v.Substring(0, $v.Length - 4) + "_FV2404"
        }
    }
}

# Turn the used range into an Excel Table ("Table1")
$rng = $ws.Range("A1:U57")
$lo = $ws.ListObjects.Add(1, $rng, [System.Reflection.Missing]::Value, 1)
$lo.Name = "Table1"
$lo.TableStyle = ""

# Freeze the header row (split below row 1)
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
